$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: fill in the xpath/selector values for the mobile profile page ---
# NB: cells are written in the exact order the author typed them (O2 before N2)
# so that the shared-string table gets the same insertion order as the source file.
$ws.Range("B2").Value = '//android.widget.Image[@resource-id="uploadedAvatar"]'
$ws.Range("C2").Value = '//android.widget.Button[@resource-id="profilePictureInput"]'
$ws.Range("D2").Value = '???'
$ws.Range("D2").Interior.Color = 49407
$ws.Range("E2").Value = '//android.widget.Button[@text=""]'
$ws.Range("F2").Value = '//android.view.View[@resource-id="gender"]'
$ws.Range("G2").Value = '//android.widget.EditText[@resource-id="first_name"]'
$ws.Range("H2").Value = '//android.widget.EditText[@resource-id="last_name"]'
$ws.Range("I2").Value = '//android.widget.EditText[@resource-id="email"]'
$ws.Range("J2").Value = '//android.widget.EditText[@resource-id="emails_copies"]'
$ws.Range("K2").Value = '//android.widget.EditText[@resource-id="phone"]'
$ws.Range("L2").Value = '//android.widget.EditText[@resource-id="mobile_phone"]'
$ws.Range("M2").Value = '//android.widget.Spinner[@resource-id="birthday"]'
$ws.Range("O2").Value = '//android.widget.CheckBox[@resource-id="hide_birthday"]'
$ws.Range("N2").Value = '//android.widget.CheckBox[@resource-id="public_contact_data"]'
$ws.Range("P2").Value = '//android.widget.EditText[@resource-id="current_password"]'
$ws.Range("Q2").Value = '//android.widget.EditText[@resource-id="new_password"]'
$ws.Range("R2").Value = '//android.widget.EditText[@resource-id="password_confirmation"]'
$ws.Range("S2").Value = '//android.widget.Button[@text="Speichern"]'

# --- Columns B:S got wider after the new, longer text was dropped in and the
# columns were auto-fit to the new content. D and D-width stay the same. ---
$ws.Columns.Item(2).ColumnWidth = 48.88671875
$ws.Columns.Item(3).ColumnWidth = 51.88671875
$ws.Columns.Item(5).ColumnWidth = 31.109375
$ws.Columns.Item(6).ColumnWidth = 38.6640625
$ws.Columns.Item(7).ColumnWidth = 46.33203125
$ws.Columns.Item(8).ColumnWidth = 46
$ws.Columns.Item(9).ColumnWidth = 42
$ws.Columns.Item(10).ColumnWidth = 49.109375
$ws.Columns.Item(11).ColumnWidth = 42.6640625
$ws.Columns.Item(12).ColumnWidth = 49.33203125
$ws.Columns.Item(13).ColumnWidth = 43.5546875
$ws.Columns.Item(14).ColumnWidth = 55.77734375
$ws.Columns.Item(15).ColumnWidth = 50.109375
$ws.Columns.Item(16).ColumnWidth = 52.44140625
$ws.Columns.Item(17).ColumnWidth = 49.88671875
$ws.Columns.Item(18).ColumnWidth = 57.21875
$ws.Columns.Item(19).ColumnWidth = 38.21875

# --- Move/realign the screenshot picture that sits behind the field map ---
# (columns are resized first so the two-cell anchor recomputes against the
# new widths, same as it would in real Excel)
$shp = $ws.Shapes.Item(1)
$shp.Left = 1.8
$shp.Top = 87.6
$shp.Width = 1710.5361417322836
$shp.Height = 887.8889763779528

# --- Scroll the sheet over to the right / change the selection, like the author did ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 13
$win.ScrollRow = 1
$ws.Range("R22").Select() | Out-Null
